$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append two new rows of data (rows 40 and 41)
$ws.Range("A40").Value = 43
$ws.Range("B40").Value = 2024
$ws.Range("A41").Value = 33
$ws.Range("B41").Value = 2024

# Update the view: scroll so row 13 is the top-left visible row,
# and select cell L28
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("L28").Select()
